$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# Update the dropdown selections in row 3 (B3:F3) to the new modifier picks.
$ws.Range("B3").Value = "All"
$ws.Range("C3").Value = "Projectile"
$ws.Range("D3").Value = "All"
$ws.Range("E3").Value = "Pierce"
$ws.Range("F3").Value = "All"

# Move the active selection to F3, matching the saved cursor position.
$ws.Range("F3").Select() | Out-Null
